$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = 232

$ws.Range("E7").Value = "Chaos Orb`nPath of Exile"

$ws.Range("H8").Value = "https://dunkbin.com/img/243.png"

$ws.Range("H10").Value = ""

$ws.Range("D13").Value = "zakzak_channel"
$ws.Range("E13").Value = 155860288
$ws.Range("F13").Value = "1009 days"

$ws.Range("D15").Value = "Galtz"
